$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper picked up two additional Yuzvendra Chahal matches (against
# Kings XI Punjab on 2020-09-24 and Delhi Capitals on 2020-10-05) that
# duplicate the stats already recorded in rows 2-3, appended here as new
# rows 4 and 5. Numeric-looking figures are kept as text (matching the
# rest of the sheet, which stores everything as strings), so the number
# format is forced to Text before the values are written.

# Row 4
$ws.Range("A4").Value = " Dubai (DSC)"
$ws.Range("B4").Value = " September 24 2020"
$ws.Range("C4").Value = "Kings XI won by 97 runs"
$ws.Range("D4").Value = "Royal Challengers Bangalore"
$ws.Range("E4").Value = "Kings XI Punjab"
$ws.Range("F4").Value = "Yuzvendra Chahal "

$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "1"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "3"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "0"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "0"
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "33.33"

# Row 5
$ws.Range("A5").Value = " Dubai (DSC)"
$ws.Range("B5").Value = " October 05 2020"
$ws.Range("C5").Value = "Capitals won by 59 runs"
$ws.Range("D5").Value = "Royal Challengers Bangalore"
$ws.Range("E5").Value = "Delhi Capitals"
$ws.Range("F5").Value = "Yuzvendra Chahal "

$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "0"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "0"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "0"
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "-"
